# "Generate Report for Archive"
# - Flip the localization status of both handed-off files from
#   "Ready for handoff" to "In Translation" on all three sheets
#   (Overview, zh-cn, de-de).
# - Narrow the "Status" column (and the two Overview columns that mirror
#   it) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column widths: ~17.22 -> ~13.41 -------------------------------------
# Excel stores column width snapped to the workbook's pixel grid, so the
# closest attainable ColumnWidth is used to land as near the target
# (13.4101845877511 "characters") as the grid allows.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe.Range("C1").ColumnWidth = 12.5
